# EDT2 version: new "adaptive" FEARFUL emotion added to the translation dictionary.
# Insert a new row right after the existing SADDER row (row 13), pushing the
# ENTER_ID..VALUE rows down by one, and fill in the new emotion's translations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 14 (row 14 "ENTER_ID" and everything below
# shifts down to make room for the new "FEARFUL" entry).
$ws.Rows.Item(14).Insert()

# key / DE / EN / RU / NL columns for the new FEARFUL row.
$ws.Range("A14").Value = "FEARFUL"
$ws.Range("B14").Value = "ängstlicher"
$ws.Range("C14").Value = "more fearful"
$ws.Range("E14").Value = "angstig"
$ws.Range("D14").Value = "более пугающей"

# The RU cell for the new row got a dedicated font in the source workbook.
$ws.Range("D14").Font.Name = "Calibri (Textkörper)"

# New custom width for column D, matching the widened RU column in the
# updated sheet.
$ws.Columns.Item(4).ColumnWidth = 51.5

# Reflect the author's final selection on the newly-added cell.
$ws.Range("D14").Select() | Out-Null
